$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add the new "CDD" abbreviation row right after the existing last data row (WBS, row 20)
$ws.Range("A21").Value = "CDD"
$ws.Range("B21").Value = "Component Design Document"

# Match the formatting used by the other abbreviation rows (copy format from row 20)
$ws.Range("A20:B20").Copy()
$ws.Range("A21:B21").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Give the new row the same row height as the rest of the table
$ws.Rows.Item(21).RowHeight = 15.75

# Reflect the new active cell selection
$ws.Range("B21").Select()
